{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertText(\"\u0418\u0437\u043c\u0435\u043d\u0438\u043b\u0430 \u0444\u0430\u0439\u043b\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$lastPara = $d.Paragraphs.Last\n$lastPara.Range.Text = \"\u0418\u0437\u043c\u0435\u043d\u0438\u043b\u0430 \u0444\u0430\u0439\u043b\"\n"}
